$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.951.03"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.749.56"
$ws.Range("E3").Value = "  +2.29%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.93"
$ws.Range("E5").Value = "  +0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.17"
$ws.Range("E6").Value = "  +0.86%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.748.40"
$ws.Range("E7").Value = "  +2.40%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  +1.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +2.48%  "
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("E12").Value = "  +0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.86"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("E14").Value = "  +1.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.376.90"
$ws.Range("E15").Value = "  +2.26%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.748.52"
$ws.Range("E16").Value = "  +2.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.989.48"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("E18").Value = "  +2.23%  "
$ws.Range("E19").Value = "  -1.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.22"
$ws.Range("E20").Value = "  +1.20%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.80"
$ws.Range("E21").Value = "  +19.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "491.80"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.725"
$ws.Range("E23").Value = "  +1.12%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000152"
$ws.Range("E24").Value = "  +7.86%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "84.75"
$ws.Range("E25").Value = "  +0.55%  "
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("E27").Value = "  +1.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.11"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E30").Value = "  +3.07%  "
$ws.Range("E31").Value = "  +4.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.01"
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.54"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.895.03"
$ws.Range("E34").Value = "  +2.42%  "
$ws.Range("E35").Value = "  +0.59%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.684.90"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("E40").Value = "  +2.10%  "
$ws.Range("E41").Value = "  +1.37%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.95"
$ws.Range("E42").Value = "  +5.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "428.85"
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.63"
$ws.Range("E44").Value = "  -0.49%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.00"
$ws.Range("E45").Value = "  +2.53%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.47"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.30"
$ws.Range("E48").Value = "  +0.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.32"
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.780.00"
$ws.Range("E50").Value = "  +2.34%  "
$ws.Range("E51").Value = "  +1.34%  "
